$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L2").Value = "[0.3809578879251784, 0.4591195577891941]"
$ws.Range("P2").Value = "[-1.3082107546480026, -1.1069475616252324]"
$ws.Range("T2").Value = "[0.3956119468871589, 0.4393171727151461]"
$ws.Range("X2").Value = 4.207087087087193
$ws.Range("Y2").Value = 4.972012012012139
